$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(3, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(4, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(5, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(6, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(7, 3.286832544864788, 1.655778082260271, 22.3905356188092, 10.19245300693656, 37.52559925287081),
    @(8, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(9, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(10, 1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302),
    @(11, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(12, 0.04271373187048222, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 4.115428400803308),
    @(13, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(14, 0.1190320826869504, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0.8034070775528621),
    @(15, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(16, 3.286832544864788, 10.34677158129881, 22.3905356188092, 10.19245300693656, 46.21659275190936),
    @(17, 0.1190320826869504, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.457851494814137),
    @(18, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(19, 0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 16.0466451790395),
    @(20, 0.1190320826869504, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.418468675747795),
    @(21, 0.2917716402565462, 3286.919754855326, 0.1494219747398047, 10.19245300693656, 3297.553401477259),
    @(22, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(23, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447),
    @(24, 3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(25, 3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(26, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(27, 1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 7.143138311642302),
    @(28, 0.1190320826869504, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.418468675747795),
    @(29, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(30, 1.455362044514542, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 12.70737954648466),
    @(31, 3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(32, 1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 3.009163075608874),
    @(33, 0.6606524410359556, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 6.348428708163715),
    @(34, 0.04271373187048222, 0.002571899574220771, 0.7527432677738641, 0.4942365360607697, 1.292265435279337),
    @(35, 3.286832544864788, 10.34677158129881, 0.7527432677738641, 10.19245300693656, 24.57880040087402),
    @(36, 3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833),
    @(37, 0.2917716402565462, 0.306821227259698, 0.7527432677738641, 10.19245300693656, 11.54378914222666),
    @(38, 0.1190320826869504, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.069511820747223),
    @(39, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(40, 0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 16.0466451790395),
    @(41, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(42, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(43, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548),
    @(44, 0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(45, 3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548),
    @(46, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(47, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(48, 0.6606524410359556, 117.745847958593, 261.3203778131603, 10.19245300693656, 389.9193312197258),
    @(49, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(50, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

Write-Host "Updated $($data.Count) rows"